# Add a new "T_VIV_COL" column to the Manzanos report, between the
# existing "T_VIV_DES" (col J) and "TOTAL_VIV" (col K) columns.
#
# This is equivalent to inserting a new column at K (shifting the old
# K/L -> L/M) and filling in the new header label in the header row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K - shifts TOTAL_VIV/TOTAL_POB right by one.
$ws.Columns.Item(11).Insert()

# New header cell for the inserted column.
$ws.Cells.Item(5, 11).Value = "T_VIV_COL"

# The insert operation clones formatting from the neighbouring column into
# row 7 (which, for every other column, is intentionally blank outside
# F:J) - clear that stray cell so the row stays F:J only, matching the
# original layout.
$ws.Cells.Item(7, 11).Clear()

# Move the selection to where the author left it after adding the column.
$ws.Range("K14").Select()
